$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1. "can start up to five services" -> "can start a number of services"
# ----------------------------------------------------------------------
$d.Content.Find.Execute(
    "can start up to five services",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "can start a number of services",
    2) | Out-Null

# ----------------------------------------------------------------------
# 2. fix typo: "before it likks off RamEater services" ->
#    "before it kills off RamEater services"
# ----------------------------------------------------------------------
$d.Content.Find.Execute(
    "before it likks off RamEater services",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "before it kills off RamEater services",
    2) | Out-Null

# ----------------------------------------------------------------------
# 3. Append new sentence after "...does kill one of the services."
# ----------------------------------------------------------------------
$d.Content.Find.Execute(
    "these notifications will disappear if the OS does kill one of the services. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "these notifications will disappear if the OS does kill one of the services. Each service will attempt to allocate the maximum amount of memory allowed, but see notes below for Android 2. ",
    2) | Out-Null

# ----------------------------------------------------------------------
# 4. New help content at the end of the document (around the _GoBack
#    bookmark): add an "Android 2" heading/paragraph, then a following
#    paragraph explaining memory allocation on older Android versions.
# ----------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bmRange = $bm.Range

# Split the bookmark's (currently empty) paragraph into two paragraphs;
# the new, first paragraph will receive the "Android 2" text and the
# bookmark stays behind in the second (original) paragraph.
$bmRange.InsertParagraphBefore()

$paras = $d.Paragraphs
$androidPara = $paras.Item($paras.Count - 1)
$androidPara.Range.InsertAfter("Android 2")

# Remove the bookmark so it can be re-anchored around the new sentence
# (re-inserting text at a collapsed bookmark always lands outside it).
$bm2 = $d.Bookmarks.Item("_GoBack")
$bm2.Delete()

$explainPara = $d.Paragraphs.Last
$explainPara.Range.InsertAfter("In Android versions earlier than Android 3 it is not possible to automatically detect the maximum amount of memory to allocate, you will need to select the amount of memory each service will allocate in the settings.")

# Re-create the _GoBack bookmark spanning the whole explanatory sentence,
# matching the bookmarkStart/bookmarkEnd wrapping from the source edit.
$explainPara2 = $d.Paragraphs.Last
$d.Bookmarks.Add("_GoBack", $explainPara2.Range) | Out-Null

Write-Output "done"
